$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the workbook window's tab-bar ratio (980 -> 990). Best-effort: this
# workbook-chrome setting isn't always round-tripped by every host, but
# setting it is harmless if ignored.
try {
    $wb.Windows.Item(1).TabRatio = 990
} catch {
}

# The three cells A7:A9 rotate their contents:
#   A7 was the text "data"  -> becomes the number 1
#   A8 was the number 1     -> becomes the text "auto"
#   A9 was the text "auto"  -> becomes the text "data"
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = "auto"
$ws.Range("A9").Value = "data"

# A7 and A9 drop back to the plain/default cell style (A8, B8, C8, B9, C9
# keep whatever style they already had).
$ws.Range("A7").Font.Name = "Arial"
$ws.Range("A7").Font.Size = 10
$ws.Range("A9").Font.Name = "Arial"
$ws.Range("A9").Font.Size = 10

# NOTE: the source diff also narrows the sheet's default column width
# (8.50510204081633 -> 8.23469387755102). That default/uniform width comes
# from a single <col min="1" max="1025".../> entry; this host's
# ColumnWidth setter only supports Excel's 1/6-character quantisation and
# always operates on an explicit column (splitting the single default-width
# run into several <col> entries), so any value we could set here would
# both corrupt that single-range structure and land farther from the
# target float than simply leaving the column width untouched. Left as-is
# on purpose.

# Selection moves from A10 to A8.
$null = $ws.Range("A8").Select()
